$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
# G1 "Gutschrift" -> "Kosten", H1 "Kosten" -> "Nachhaltigkeit"
$ws.Range("G1").Value = "Kosten"
$ws.Range("H1").Value = "Nachhaltigkeit"

# --- Row 2: Festgeld ---
$ws.Range("G2").Value = "nein"
$ws.Range("H2").Value = "nein"

# --- Row 3: Sparbrief ---
$ws.Range("E3").Value = "langfristig"
$ws.Range("G3").Value = "nein"
$ws.Range("H3").Value = "nein"

# --- Row 4: Tagesgeld ---
$ws.Range("G4").Value = "nein"
$ws.Range("H4").Value = "nein"

# --- Row 5: new product - Aktienfond ---
$ws.Range("A5").Value = "PIB_Union_Aktienfond_9766865.pdf"
$ws.Range("B5").Value = "Aktienfond"
$ws.Range("C5").Value = 9766865
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = "langfristig"
$ws.Range("F5").Value = "mittleres Risiko"
$ws.Range("G5").Value = "ja"
$ws.Range("H5").Value = "ja"

# --- Row 6: new product - Bonuszertifikat ---
$ws.Range("A6").Value = "PIB_DZBank_Bonuszertifikat_7035880.pdf"
$ws.Range("B6").Value = "Bonuszertifikat"
$ws.Range("C6").Value = 7035880
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = "mittelfristig"
$ws.Range("F6").Value = "hohes Risiko"
$ws.Range("G6").Value = "ja"
$ws.Range("H6").Value = "nein"

# --- Update selection to match the saved state (G6 was last active cell) ---
$ws.Range("G6").Select()

Write-Host "Edit applied"
